$wb = $excel.ActiveWorkbook

# Insert the new "UpdateActivity" worksheet right after the "Users" sheet
$usersSheet = $wb.Worksheets.Item("Users")
$newSheet = $wb.Worksheets.Add($null, $usersSheet)
$newSheet.Name = "UpdateActivity"

# Header row
$newSheet.Range("A1").Value = "Subject"
$newSheet.Range("B1").Value = "IndustryGroup"
$newSheet.Range("C1").Value = "ProductType"
$newSheet.Range("D1").Value = "Description"
$newSheet.Range("E1").Value = "MeetingNotes"
$newSheet.Range("F1").Value = "ExtAttendee"
$newSheet.Range("G1").Value = "HLAttendee"
$newSheet.Range("A1:G1").Font.Bold = $true
$newSheet.Range("A1:G1").HorizontalAlignment = -4108

# Data row (values are assigned in this particular order so that newly
# introduced shared strings land at the same table positions as the target)
$newSheet.Range("B2").Value = "FIG - Financial Institutions"
$newSheet.Range("C2").Value = "Advisory"
$newSheet.Range("E2").Value = "Updated Notes"
$newSheet.Range("G2").Value = "Amanda Donovan"
$newSheet.Range("A2").Value = "Updated By Delegate"
$newSheet.Range("D2").Value = "Updated By Delegate Description"
$newSheet.Range("F2").Value = "Test James"

# (Column widths below are the closest values this engine's ColumnWidth
# setter can reproduce toward the authored best-fit widths of 27.88671875,
# 25.21875, 14.33203125, 43, 14.5546875, 11.5546875 and 17.6640625.)
$newSheet.Columns.Item(1).ColumnWidth = 27
$newSheet.Columns.Item(2).ColumnWidth = 24.333333333333332
$newSheet.Columns.Item(3).ColumnWidth = 13.5
$newSheet.Columns.Item(4).ColumnWidth = 42.166666666666664
$newSheet.Columns.Item(5).ColumnWidth = 13.666666666666666
$newSheet.Columns.Item(6).ColumnWidth = 10.666666666666666
$newSheet.Columns.Item(7).ColumnWidth = 16.833333333333332

# Update selection on the "Users" sheet (it is no longer the active tab)
$usersSheet.Range("H18").Select() | Out-Null

# "UpdateActivity" ends up being the active/selected sheet
$newSheet.Activate()
$newSheet.Range("D8").Select() | Out-Null
